# Apply updated cryptocurrency price/volume data (and two row swaps)
# to the worksheet, forcing Text storage (matching the original inline-string
# cell type) instead of letting General-format auto-convert numeric-looking
# strings (e.g. "0.9993", "242.79") into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" '29.162.37'
Set-TextCell "E2" '  +0.34%  '

Set-TextCell "D3" '1.829.54'
Set-TextCell "E3" '  -0.23%  '

Set-TextCell "D4" '0.9993'
Set-TextCell "E4" '  -0.04%  '

Set-TextCell "D5" '242.79'
Set-TextCell "E5" '  +0.08%  '

Set-TextCell "D6" '0.6162'
Set-TextCell "E6" '  -0.33%  '

Set-TextCell "D7" '1.001'
Set-TextCell "E7" '  -0.09%  '

Set-TextCell "D8" '0.07339'
Set-TextCell "E8" '  -1.64%  '

Set-TextCell "D9" '0.2911'
Set-TextCell "E9" '  -0.51%  '

Set-TextCell "D10" '23.16'
Set-TextCell "E10" '  +0.38%  '

Set-TextCell "D11" '0.07636'
Set-TextCell "E11" '  -0.56%  '

Set-TextCell "D12" '1.829.73'
Set-TextCell "E12" '  -0.09%  '

Set-TextCell "E13" '  -0.55%  '

Set-TextCell "D14" '0.6701'
Set-TextCell "E14" '  -0.24%  '

Set-TextCell "D15" '82.32'
Set-TextCell "E15" '  -0.34%  '

Set-TextCell "E16" '  -2.30%  '

Set-TextCell "D17" '5.838'
Set-TextCell "E17" '  -1.16%  '

Set-TextCell "D18" '29.147.20'
Set-TextCell "E18" '  +0.33%  '

Set-TextCell "D19" '2.071.21'
Set-TextCell "E19" '  -0.63%  '

Set-TextCell "D20" '236.08'
Set-TextCell "E20" '  +1.53%  '

Set-TextCell "D21" '12.47'
Set-TextCell "E21" '  -1.70%  '

Set-TextCell "E22" '  -0.14%  '

Set-TextCell "D23" '7.355'
Set-TextCell "E23" '  +2.08%  '

Set-TextCell "D24" '1.002'
Set-TextCell "E24" '  +0.05%  '

Set-TextCell "D25" '158.61'
Set-TextCell "E25" '  -0.44%  '

Set-TextCell "D26" '0.1387'
Set-TextCell "E26" '  -2.22%  '

Set-TextCell "D27" '8.508'
Set-TextCell "E27" '  +0.36%  '

Set-TextCell "D28" '17.60'
Set-TextCell "E28" '  -1.05%  '

Set-TextCell "D29" '1.487'
Set-TextCell "E29" '  -0.97%  '

Set-TextCell "D30" '0.05864'
Set-TextCell "E30" '  +5.93%  '

Set-TextCell "E31" '  +1.42%  '

Set-TextCell "E32" '  -0.61%  '

Set-TextCell "E33" '  -1.78%  '

Set-TextCell "D34" '1.855'
Set-TextCell "E34" '  +1.09%  '

Set-TextCell "B35" 'ARBITRUM'
Set-TextCell "C35" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell "D35" '1.136'
Set-TextCell "E35" '  -0.27%  '

Set-TextCell "B36" 'ImmutableX'
Set-TextCell "C36" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell "D36" '0.7235'
Set-TextCell "E36" '  -1.98%  '

Set-TextCell "D37" '2.611'
Set-TextCell "E37" '  -1.93%  '

Set-TextCell "D38" '2.864'
Set-TextCell "E38" '  +3.25%  '

Set-TextCell "D39" '1.225.72'
Set-TextCell "E39" '  +1.44%  '

Set-TextCell "E40" '  -1.13%  '

Set-TextCell "D41" '6.188'
Set-TextCell "E41" '  -4.19%  '

Set-TextCell "D42" '0.9034'

Set-TextCell "E43" '  +0.04%  '

Set-TextCell "D44" '101.80'

Set-TextCell "D45" '1.984.13'
Set-TextCell "E45" '  +0.30%  '

Set-TextCell "D46" '65.64'
Set-TextCell "E46" '  +0.32%  '

Set-TextCell "D47" '0.5046'
Set-TextCell "E47" '  -0.86%  '

Set-TextCell "B48" 'EnergySwap'
Set-TextCell "C48" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell "D48" '9.189'
Set-TextCell "E48" '  +0.60%  '

Set-TextCell "B49" 'TheSandbox'
Set-TextCell "C49" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell "D49" '0.4043'
Set-TextCell "E49" '  -0.64%  '

Set-TextCell "E50" '  -3.19%  '

Set-TextCell "D51" '0.1142'
Set-TextCell "E51" '  +3.52%  '
